$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3099390012751145
$ws.Range("J2").Value = 0.3099390012751145
$ws.Range("M2").Value = 68.63737500000001
$ws.Range("N2").Value = 205.912125
$ws.Range("O2").Value = 0.5415701538216162
$ws.Range("P2").Value = 0.5415701538216162
$ws.Range("Q2").Value = 8.230536427500001
$ws.Range("R2").Value = 74.0748278475
$ws.Range("S2").Value = 0.1678537125958819
$ws.Range("T2").Value = 0.1678537125958819
$ws.Range("I3").Value = 0.3099390012751145
$ws.Range("J3").Value = 0.3099390012751145
$ws.Range("O3").Value = 0.08718851262838957
$ws.Range("P3").Value = 0.08718851262838957
$ws.Range("S3").Value = 0.02702312052670577
$ws.Range("T3").Value = 0.02702312052670577
$ws.Range("I4").Value = 0.3099390012751145
$ws.Range("J4").Value = 0.3099390012751145
$ws.Range("M4").Value = 16.21089566666667
$ws.Range("N4").Value = 48.632687
$ws.Range("O4").Value = 0.1279089892319285
$ws.Range("P4").Value = 0.1279089892319285
$ws.Range("Q4").Value = 1.943902535708889
$ws.Range("R4").Value = 17.49512282138
$ws.Range("S4").Value = 0.0396439843766533
$ws.Range("T4").Value = 0.0396439843766533
$ws.Range("I5").Value = 0.3099390012751145
$ws.Range("J5").Value = 0.3099390012751145
$ws.Range("M5").Value = 20.32546233333333
$ws.Range("N5").Value = 60.976387
$ws.Range("O5").Value = 0.1603741949973873
$ws.Range("P5").Value = 0.1603741949973873
$ws.Range("Q5").Value = 2.437293939931111
$ws.Range("R5").Value = 21.93564545938
$ws.Range("S5").Value = 0.0497062178277907
$ws.Range("T5").Value = 0.0497062178277907
$ws.Range("I6").Value = 0.3099390012751145
$ws.Range("J6").Value = 0.3099390012751145
$ws.Range("M6").Value = 10.513928
$ws.Range("N6").Value = 31.541784
$ws.Range("O6").Value = 0.08295814932067838
$ws.Range("P6").Value = 0.08295814932067838
$ws.Range("Q6").Value = 1.260760152906667
$ws.Range("R6").Value = 11.34684137616
$ws.Range("S6").Value = 0.02571196594808288
$ws.Range("T6").Value = 0.02571196594808288
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2669800000000001
$ws.Range("H7").Value = 0.8009400000000001
$ws.Range("I7").Value = 0.6900609987248855
$ws.Range("J7").Value = 0.6900609987248854
$ws.Range("M7").Value = 68.63737500000001
$ws.Range("N7").Value = 205.912125
$ws.Range("O7").Value = 0.5415701538216162
$ws.Range("P7").Value = 0.5415701538216162
$ws.Range("Q7").Value = 18.3248063775
$ws.Range("R7").Value = 164.9232573975
$ws.Range("S7").Value = 0.3737164412257343
$ws.Range("T7").Value = 0.3737164412257343
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2669800000000001
$ws.Range("H8").Value = 0.8009400000000001
$ws.Range("I8").Value = 0.6900609987248855
$ws.Range("J8").Value = 0.6900609987248854
$ws.Range("O8").Value = 0.08718851262838957
$ws.Range("P8").Value = 0.08718851262838957
$ws.Range("Q8").Value = 2.950148934506667
$ws.Range("R8").Value = 26.55134041056001
$ws.Range("S8").Value = 0.0601653921016838
$ws.Range("T8").Value = 0.06016539210168378
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2669800000000001
$ws.Range("H9").Value = 0.8009400000000001
$ws.Range("I9").Value = 0.6900609987248855
$ws.Range("J9").Value = 0.6900609987248854
$ws.Range("M9").Value = 16.21089566666667
$ws.Range("N9").Value = 48.632687
$ws.Range("O9").Value = 0.1279089892319285
$ws.Range("P9").Value = 0.1279089892319285
$ws.Range("Q9").Value = 4.327984925086668
$ws.Range("R9").Value = 38.95186432578001
$ws.Range("S9").Value = 0.08826500485527521
$ws.Range("T9").Value = 0.08826500485527519
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2669800000000001
$ws.Range("H10").Value = 0.8009400000000001
$ws.Range("I10").Value = 0.6900609987248855
$ws.Range("J10").Value = 0.6900609987248854
$ws.Range("M10").Value = 20.32546233333333
$ws.Range("N10").Value = 60.976387
$ws.Range("O10").Value = 0.1603741949973873
$ws.Range("P10").Value = 0.1603741949973873
$ws.Range("Q10").Value = 5.426491933753335
$ws.Range("R10").Value = 48.83842740378001
$ws.Range("S10").Value = 0.1106679771695966
$ws.Range("T10").Value = 0.1106679771695966
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2669800000000001
$ws.Range("H11").Value = 0.8009400000000001
$ws.Range("I11").Value = 0.6900609987248855
$ws.Range("J11").Value = 0.6900609987248854
$ws.Range("M11").Value = 10.513928
$ws.Range("N11").Value = 31.541784
$ws.Range("O11").Value = 0.08295814932067838
$ws.Range("P11").Value = 0.08295814932067838
$ws.Range("Q11").Value = 2.80700849744
$ws.Range("R11").Value = 25.26307647696
$ws.Range("S11").Value = 0.0572461833725955
$ws.Range("T11").Value = 0.0572461833725955
